$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.439.01'
$ws.Range("E2").Value = '  -0.38%  '

$ws.Range("D3").Value = '1.799.82'
$ws.Range("E3").Value = '  -0.73%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.40'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.602'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.04'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +6.45%  '

$ws.Range("E9").Value = '  -4.58%  '

$ws.Range("E10").Value = '  -4.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.76%  '

$ws.Range("D12").Value = '2.059.36'
$ws.Range("E12").Value = '  -0.68%  '

$ws.Range("D13").Value = '1.799.98'
$ws.Range("E13").Value = '  -1.37%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.87'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.51%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.628'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -4.21%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '34.415.93'
$ws.Range("E16").Value = '  -0.39%  '

$ws.Range("E17").Value = '  -3.40%  '

$ws.Range("E18").Value = '  -3.48%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '238.86'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.38%  '

$ws.Range("E20").Value = '  -4.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.06'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.09%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.07'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.16'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.92%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.54'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.55%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.65'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.92%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.51'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.120'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.46%  '

$ws.Range("E29").Value = '  -0.02%  '

$ws.Range("E30").Value = '  -1.82%  '

$ws.Range("E31").Value = '  -3.71%  '

$ws.Range("E32").Value = '  -3.79%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.82'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.80'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.30%  '

$ws.Range("D35").Value = '1.306.34'
$ws.Range("E35").Value = '  -6.96%  '

$ws.Range("E36").Value = '  -5.40%  '

$ws.Range("E37").Value = '  -1.77%  '

$ws.Range("E38").Value = '  -3.02%  '

$ws.Range("E39").Value = '  -6.59%  '

$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.22'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.35%  '

$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.44'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.64'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.76%  '

$ws.Range("E43").Value = '  -1.82%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.942'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.95%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.04'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0516'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +4.05%  '

$ws.Range("D47").Value = '1.961.23'
$ws.Range("E47").Value = '  -0.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.70'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.89%  '

$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.76'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.77%  '

$ws.Range("E51").Value = '  -0.70%  '
